# chore: update Sheets via scheduled runner
#
# Applies the numeric corrections captured in the commit diff to the
# four affected worksheets (BSM, CRP, CUL, LTW, WVR) of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# BSM: two rows get their "found" profit columns bumped up, and a long run
# of rows that previously priced out to a flat loss lose their cached
# H..N figures entirely (the cells go back to being empty / not present).
# ---------------------------------------------------------------------------
$bsm = $wb.Worksheets.Item("BSM")

$bsm.Range("H76").Value2 = 11562.6
$bsm.Range("J76").Value2 = 11562.6
$bsm.Range("L76").Value2 = 11562.6
$bsm.Range("N76").Value2 = -12192.6

$bsm.Range("H79").Value2 = 11562.6
$bsm.Range("J79").Value2 = 11562.6
$bsm.Range("L79").Value2 = 11562.6
$bsm.Range("N79").Value2 = -13746.6

$bsmClearRows = @(117,118,119,120,122,123,124,125,126,127,128,129,130,131,132,133,134,135,137,138,139,140,141)
foreach ($r in $bsmClearRows) {
    $bsm.Range("H" + $r + ":N" + $r).ClearContents()
}

# ---------------------------------------------------------------------------
# CRP: two previously-empty rows now carry a filled-in 22222 cost basis.
# ---------------------------------------------------------------------------
$crp = $wb.Worksheets.Item("CRP")

$crp.Range("H82").Value2 = 22222
$crp.Range("J82").Value2 = 22222
$crp.Range("L82").Value2 = 22222
$crp.Range("N82").Value2 = -22944

$crp.Range("H85").Value2 = 22222
$crp.Range("J85").Value2 = 22222
$crp.Range("L85").Value2 = 22222
$crp.Range("N85").Value2 = -24718

# ---------------------------------------------------------------------------
# CUL: row 113's blended-cost figures are recalculated.
# ---------------------------------------------------------------------------
$cul = $wb.Worksheets.Item("CUL")

$cul.Range("H113").Value2 = 1091.3667
$cul.Range("I113").Value2 = 868.7143
$cul.Range("J113").Value2 = 1610.8889
$cul.Range("K113").Value2 = 2606.1429
$cul.Range("L113").Value2 = 4832.6667
$cul.Range("M113").Value2 = -436.1428999999998
$cul.Range("N113").Value2 = -9172.6667

# ---------------------------------------------------------------------------
# LTW: same kind of cleanup as BSM above - a run of rows loses its cached
# H..N figures entirely.
# ---------------------------------------------------------------------------
$ltw = $wb.Worksheets.Item("LTW")

$ltwClearRows = @(124,125,127,128,129,130,131,132,133,134,135,136,137,138,139,140,141)
foreach ($r in $ltwClearRows) {
    $ltw.Range("H" + $r + ":N" + $r).ClearContents()
}

# ---------------------------------------------------------------------------
# WVR: mirrors the CRP change - two previously-empty rows now carry a
# filled-in 23000 cost basis.
# ---------------------------------------------------------------------------
$wvr = $wb.Worksheets.Item("WVR")

$wvr.Range("H82").Value2 = 23000
$wvr.Range("J82").Value2 = 23000
$wvr.Range("L82").Value2 = 23000
$wvr.Range("N82").Value2 = -23766

$wvr.Range("H85").Value2 = 23000
$wvr.Range("J85").Value2 = 23000
$wvr.Range("L85").Value2 = 23000
